$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 999.5
$ws.Range("I43").Value = 999
$ws.Range("K43").Value = 999
$ws.Range("M43").Value = -930

$ws.Range("H62").Value = 8212.143
$ws.Range("I62").Value = 8214.333000000001
$ws.Range("J62").Value = 8199
$ws.Range("K62").Value = 8214.333000000001
$ws.Range("L62").Value = 8199
$ws.Range("M62").Value = -7590.333000000001
$ws.Range("N62").Value = -9447

$ws.Range("H65").Value = 8212.143
$ws.Range("I65").Value = 8214.333000000001
$ws.Range("J65").Value = 8199
$ws.Range("K65").Value = 41071.665
$ws.Range("L65").Value = 40995
$ws.Range("M65").Value = -37951.665
$ws.Range("N65").Value = -47235

$ws.Range("H94").Value = 1512
$ws.Range("I94").Value = 1074
$ws.Range("K94").Value = 1074
$ws.Range("M94").Value = -623

$ws.Range("H96").Value = 380.2857
$ws.Range("I96").Value = 275.1
$ws.Range("J96").Value = 643.25
$ws.Range("K96").Value = 825.3000000000001
$ws.Range("L96").Value = 1929.75
$ws.Range("M96").Value = 547.6999999999999
$ws.Range("N96").Value = -4675.75

$ws.Range("H103").Value = 1149
$ws.Range("I103").Value = 1099.5
$ws.Range("J103").Value = 1198.5
$ws.Range("K103").Value = 3298.5
$ws.Range("L103").Value = 3595.5
$ws.Range("M103").Value = -2712.5
$ws.Range("N103").Value = -4767.5

$ws.Range("H107").Value = 5134.5713
$ws.Range("I107").Value = 4323.6665
$ws.Range("K107").Value = 4323.6665
$ws.Range("M107").Value = -2403.6665

$ws.Range("H111").Value = 602.5
$ws.Range("I111").Value = 602.5
$ws.Range("K111").Value = 1807.5
$ws.Range("M111").Value = 1259.5

$ws.Range("H137").Value = 1916
$ws.Range("I137").Value = 1735.8334
$ws.Range("J137").Value = 2997
$ws.Range("K137").Value = 5207.5002
$ws.Range("L137").Value = 8991
$ws.Range("M137").Value = -2657.5002
$ws.Range("N137").Value = -14091

$ws.Range("H141").Value = 1698.3334
$ws.Range("I141").Value = 1698.3334
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5095.0002
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 84.9997999999996
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6644.4346
$ws.Range("I32").Value = 6644.4346
$ws.Range("K32").Value = 6644.4346
$ws.Range("M32").Value = -6357.4346

$ws.Range("H61").Value = 1188.4286
$ws.Range("I61").Value = 1219.8334
$ws.Range("K61").Value = 1219.8334
$ws.Range("M61").Value = -1007.8334

$ws.Range("H131").Value = 106666.336
$ws.Range("J131").Value = 106666.336
$ws.Range("L131").Value = 106666.336
$ws.Range("N131").Value = -116746.336

$ws.Range("H132").Value = 5856
$ws.Range("I132").Value = 6999
$ws.Range("K132").Value = 20997
$ws.Range("M132").Value = -18467

$ws.Range("H136").Value = 1188.4286
$ws.Range("I136").Value = 1219.8334
$ws.Range("K136").Value = 3659.5002
$ws.Range("M136").Value = -1109.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 900
$ws.Range("I11").Value = 800
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 800
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = -660
$ws.Range("N11").Value = -1280

$ws.Range("H16").Value = 3333
$ws.Range("J16").Value = 3333
$ws.Range("L16").Value = 3333
$ws.Range("N16").Value = -3907

$ws.Range("H31").Value = 3031.4285
$ws.Range("I31").Value = 2804.2
$ws.Range("J31").Value = 3599.5
$ws.Range("K31").Value = 2804.2
$ws.Range("L31").Value = 3599.5
$ws.Range("M31").Value = -2509.2
$ws.Range("N31").Value = -4189.5

$ws.Range("H34").Value = 3031.4285
$ws.Range("I34").Value = 2804.2
$ws.Range("J34").Value = 3599.5
$ws.Range("K34").Value = 2804.2
$ws.Range("L34").Value = 3599.5
$ws.Range("M34").Value = -2602.2
$ws.Range("N34").Value = -4003.5

$ws.Range("H107").Value = 1791.6154
$ws.Range("I107").Value = 2210.7778
$ws.Range("K107").Value = 2210.7778
$ws.Range("M107").Value = -290.7777999999998

$ws.Range("H113").Value = 3333
$ws.Range("J113").Value = 3333
$ws.Range("L113").Value = 3333
$ws.Range("N113").Value = -7673

$ws.Range("H141").Value = 56332.332
$ws.Range("J141").Value = 56332.332
$ws.Range("L141").Value = 56332.332
$ws.Range("N141").Value = -66692.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 5255.5
$ws.Range("I63").Value = 511
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 1533
$ws.Range("L63").Value = 30000
$ws.Range("M63").Value = -784
$ws.Range("N63").Value = -31498

$ws.Range("H64").Value = 233
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 233
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 699
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -1239

$ws.Range("H66").Value = 5255.5
$ws.Range("I66").Value = 511
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 4599
$ws.Range("L66").Value = 90000
$ws.Range("M66").Value = -855
$ws.Range("N66").Value = -97488

$ws.Range("H67").Value = 233
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 233
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 699
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -2571

$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H87").Value = 9249.75
$ws.Range("I87").Value = 6999.6665
$ws.Range("J87").Value = 16000
$ws.Range("K87").Value = 20998.9995
$ws.Range("L87").Value = 48000
$ws.Range("M87").Value = -19750.9995
$ws.Range("N87").Value = -50496

$ws.Range("H90").Value = 9249.75
$ws.Range("I90").Value = 6999.6665
$ws.Range("J90").Value = 16000
$ws.Range("K90").Value = 62996.9985
$ws.Range("L90").Value = 144000
$ws.Range("M90").Value = -56756.9985
$ws.Range("N90").Value = -156480

$ws.Range("H98").Value = 1677.8572
$ws.Range("J98").Value = 2125.6667
$ws.Range("L98").Value = 6377.000100000001
$ws.Range("N98").Value = -9373.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5977.75
$ws.Range("I102").Value = 5977.75
$ws.Range("K102").Value = 5977.75
$ws.Range("M102").Value = -4355.75

$ws.Range("H122").Value = 3999.4
$ws.Range("I122").Value = 3999.4
$ws.Range("K122").Value = 11998.2
$ws.Range("M122").Value = -9548.200000000001

$ws.Range("H132").Value = 3189.625
$ws.Range("J132").Value = 3999
$ws.Range("L132").Value = 11997
$ws.Range("N132").Value = -17057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5519.6665
$ws.Range("I22").Value = 3709.5
$ws.Range("K22").Value = 3709.5
$ws.Range("M22").Value = -3414.5

$ws.Range("H27").Value = 5519.6665
$ws.Range("I27").Value = 3709.5
$ws.Range("K27").Value = 3709.5
$ws.Range("M27").Value = -3602.5

$ws.Range("H40").Value = 4163.5625
$ws.Range("I40").Value = 4111.4165
$ws.Range("K40").Value = 4111.4165
$ws.Range("M40").Value = -3975.4165

$ws.Range("H122").Value = 5907.1
$ws.Range("I122").Value = 5907.1
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 17721.3
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15271.3
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 5451.2144
$ws.Range("I132").Value = 4085.5
$ws.Range("K132").Value = 12256.5
$ws.Range("M132").Value = -9726.5

$ws.Range("H136").Value = 3173.4614
$ws.Range("I136").Value = 3023.182
$ws.Range("K136").Value = 9069.545999999998
$ws.Range("M136").Value = -6519.545999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4149.5454
$ws.Range("I62").Value = 2925
$ws.Range("J62").Value = 4849.2856
$ws.Range("K62").Value = 2925
$ws.Range("L62").Value = 4849.2856
$ws.Range("M62").Value = -2301
$ws.Range("N62").Value = -6097.2856

$ws.Range("H65").Value = 4149.5454
$ws.Range("I65").Value = 2925
$ws.Range("J65").Value = 4849.2856
$ws.Range("K65").Value = 14625
$ws.Range("L65").Value = 24246.428
$ws.Range("M65").Value = -11505
$ws.Range("N65").Value = -30486.428
